$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "97.057.20"
Set-TextValue "E2" "  -0.32%  "
Set-TextValue "D3" "3.694.50"
Set-TextValue "E3" "  +3.14%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "240.25"
Set-TextValue "E5" "  -0.43%  "
Set-TextValue "D6" "1.89"
Set-TextValue "E6" "  +9.85%  "
Set-TextValue "D7" "655.10"
Set-TextValue "E7" "  -0.07%  "
Set-TextValue "E8" "  -1.42%  "
Set-TextValue "E9" "  +3.53%  "
Set-TextValue "D10" "1.00"
Set-TextValue "E10" "  +0.07%  "
Set-TextValue "D11" "3.692.86"
Set-TextValue "E11" "  +3.20%  "
Set-TextValue "D12" "45.59"
Set-TextValue "E12" "  +2.66%  "
Set-TextValue "E13" "  +1.03%  "
Set-TextValue "D14" "6.88"
Set-TextValue "E14" "  +7.12%  "
Set-TextValue "D15" "4.377.05"
Set-TextValue "E15" "  +3.06%  "
Set-TextValue "E16" "  +2.96%  "
Set-TextValue "D17" "96.787.46"
Set-TextValue "E17" "  -0.32%  "
Set-TextValue "D18" "9.09"
Set-TextValue "E18" "  +4.67%  "
Set-TextValue "D19" "3.697.28"
Set-TextValue "E19" "  +3.39%  "
Set-TextValue "D20" "19.25"
Set-TextValue "E20" "  +6.65%  "
Set-TextValue "D21" "12.82"
Set-TextValue "E21" "  +1.04%  "
Set-TextValue "E22" "  -0.60%  "
Set-TextValue "D23" "531.33"
Set-TextValue "E23" "  +3.11%  "
Set-TextValue "E24" "  +0.34%  "
Set-TextValue "D25" "7.13"
Set-TextValue "E25" "  +3.04%  "
Set-TextValue "E26" "  -1.09%  "
Set-TextValue "D27" "102.71"
Set-TextValue "E27" "  +0.78%  "
Set-TextValue "D28" "13.47"
Set-TextValue "E28" "  +2.71%  "
Set-TextValue "E29" "  -0.68%  "
Set-TextValue "D30" "12.54"
Set-TextValue "E30" "  +4.28%  "
Set-TextValue "D31" "3.06"
Set-TextValue "E31" "  +1.76%  "
Set-TextValue "D32" "1.00"
Set-TextValue "E32" "  +0.36%  "
Set-TextValue "D33" "1.89"
Set-TextValue "E33" "  +14.98%  "
Set-TextValue "E34" "  +0.80%  "
Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  +0.24%  "
Set-TextValue "D36" "32.76"
Set-TextValue "E36" "  +2.57%  "
Set-TextValue "D37" "654.92"
Set-TextValue "E37" "  +6.59%  "
Set-TextValue "D38" "0.608"
Set-TextValue "E38" "  +7.03%  "
Set-TextValue "D39" "9.07"
Set-TextValue "E39" "  +3.24%  "
Set-TextValue "D40" "7.00"
Set-TextValue "E40" "  +16.36%  "
Set-TextValue "D41" "0.163"
Set-TextValue "E41" "  +5.64%  "
Set-TextValue "E42" "  +2.97%  "
Set-TextValue "D43" "0.968"
Set-TextValue "E43" "  +4.44%  "
Set-TextValue "D44" "38.13"
Set-TextValue "E44" "  +15.91%  "
Set-TextValue "E45" "  +0.04%  "
Set-TextValue "E46" "  +7.55%  "
Set-TextValue "E47" "  +4.38%  "
Set-TextValue "E48" "  +0.29%  "
Set-TextValue "D49" "23.63"
Set-TextValue "E49" "  +0.05%  "
Set-TextValue "E50" "  +2.58%  "
Set-TextValue "D51" "3.63"
Set-TextValue "E51" "  +3.80%  "
